# Apply the commit's changes to the workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update URL, Version, Date, Publisher ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/process-client-id"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Sheet "Elements": update the Extension.url Fixed Value and clear the
#     Extension row's Constraint(s) cell ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/process-client-id"
$elements.Range("AI2").Value = ""
